$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell F1 ("time_taken"), matching the style used by the
# other header cells (B1:E1) - bold, centered, top-aligned, thin border.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the time_taken values for each data row (2-36) as text, matching
# the microsecond-precision timestamps recorded for this panel export.
$timestamps = @(
    "2021-10-05 10:52:27.299039",
    "2021-10-05 10:52:27.299053",
    "2021-10-05 10:52:27.299057",
    "2021-10-05 10:52:27.299060",
    "2021-10-05 10:52:27.299064",
    "2021-10-05 10:52:27.299067",
    "2021-10-05 10:52:27.299070",
    "2021-10-05 10:52:27.299073",
    "2021-10-05 10:52:27.299077",
    "2021-10-05 10:52:27.299080",
    "2021-10-05 10:52:27.299083",
    "2021-10-05 10:52:27.299089",
    "2021-10-05 10:52:27.299094",
    "2021-10-05 10:52:27.299099",
    "2021-10-05 10:52:27.299104",
    "2021-10-05 10:52:27.299108",
    "2021-10-05 10:52:27.299112",
    "2021-10-05 10:52:27.299115",
    "2021-10-05 10:52:27.299118",
    "2021-10-05 10:52:27.299121",
    "2021-10-05 10:52:27.299125",
    "2021-10-05 10:52:27.299128",
    "2021-10-05 10:52:27.299131",
    "2021-10-05 10:52:27.299134",
    "2021-10-05 10:52:27.299137",
    "2021-10-05 10:52:27.299141",
    "2021-10-05 10:52:27.299144",
    "2021-10-05 10:52:27.299147",
    "2021-10-05 10:52:27.299150",
    "2021-10-05 10:52:27.299153",
    "2021-10-05 10:52:27.299156",
    "2021-10-05 10:52:27.299160",
    "2021-10-05 10:52:27.299163",
    "2021-10-05 10:52:27.299167",
    "2021-10-05 10:52:27.299170"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
